$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix error of turm classes being only one:
# "Desenho Técnico" was shown in the wrong day/row combination, and
# "EAP" classes that no longer occur need to be cleared.

$ws.Range("F3").Value = "-"
$ws.Range("D4").Value = "Desenho Técnico"
$ws.Range("C6").Value = "Desenho Técnico"
$ws.Range("D6").Value = "-"
$ws.Range("F7").Value = "-"
$ws.Range("F13").Value = "-"
